# Generate Report for Handoff
#
# "b.md" finished a handoff cycle: its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff" on every
# sheet, a new (newer-versioned) handoff file/timestamp is recorded for
# both locales, the "Content Duplicate" flag drops to False, and an
# Error Detail message explaining the stale handback version is filled
# in. Column P ("Error Detail") is widened to fit the new text.

$wb = $excel.ActiveWorkbook

# Helper: write a literal-text value into a cell even when the text
# looks like a boolean ("True"/"False") so Excel doesn't silently
# coerce it to a Boolean cell. The leading apostrophe forces "stored as
# text"; resetting .Style afterwards drops the quote-prefix flag again
# so the cell's style stays plain/default, matching a normal text cell.
function Set-TextValue($range, [string]$text) {
    if ($text -eq "True" -or $text -eq "False" -or $text -eq "TRUE" -or $text -eq "FALSE") {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

# ---------------------------------------------------------------------
# Overview sheet: roll up b.md's new status + latest handoff timestamp
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
Set-TextValue $overview.Range("E3") "Ready for handoff"
Set-TextValue $overview.Range("F3") "Ready for handoff"
Set-TextValue $overview.Range("G3") "2016-09-01 04:42:15"

# ---------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
Set-TextValue $zhcn.Range("C3") "Ready for handoff"
Set-TextValue $zhcn.Range("F3") "False"
Set-TextValue $zhcn.Range("G3") "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
Set-TextValue $zhcn.Range("H3") "2016-09-01 04:42:11"
Set-TextValue $zhcn.Range("P3") "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae9495e52c78566d65e16f09790c8c92c691dbc5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb68e79ae7e889b2dbf65a9add8ec31b2f7b478f/e2e/b.md."

# Widen the Error Detail column now that it holds a long message.
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
Set-TextValue $dede.Range("C3") "Ready for handoff"
Set-TextValue $dede.Range("F3") "False"
Set-TextValue $dede.Range("G3") "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
Set-TextValue $dede.Range("H3") "2016-09-01 04:42:15"
Set-TextValue $dede.Range("P3") "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae9495e52c78566d65e16f09790c8c92c691dbc5/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb68e79ae7e889b2dbf65a9add8ec31b2f7b478f/e2e/b.md."

# Widen the Error Detail column now that it holds a long message.
$dede.Columns.Item(16).ColumnWidth = 39.17
